$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Albuns"

$albuns = @(
  @("(Equals)", 0),
  @("(Plus)", 1),
  @("× (Multiply)", 7),
  @("÷ (Divide)", 3),
  @("Don’t - EP", 0),
  @("Loose Change - EP", 0),
  @("No. 5 Collaborations Project - EP", 0),
  @("No.6 Collaborations Project", 0),
  @("Songs I Wrote With Amy - EP", 0),
  @("You Need Me - EP", 0)
)

$ws2.Range("A1").Value = "Album"
$ws2.Range("B1").Value = "Prêmios"

$r = 2
foreach ($pair in $albuns) {
  $ws2.Cells.Item($r, 1).Value = $pair[0]
  $ws2.Cells.Item($r, 2).Value = $pair[1]
  $r = $r + 1
}

$ws1.Activate()
